$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3024.1924
$ws.Range("I64").Value = 2821.111
$ws.Range("J64").Value = 3066.6978
$ws.Range("K64").Value = 2821.111
$ws.Range("L64").Value = 3066.6978
$ws.Range("M64").Value = -2573.111
$ws.Range("N64").Value = -3562.6978
# Row 67
$ws.Range("H67").Value = 3024.1924
$ws.Range("I67").Value = 2821.111
$ws.Range("J67").Value = 3066.6978
$ws.Range("K67").Value = 2821.111
$ws.Range("L67").Value = 3066.6978
$ws.Range("M67").Value = -1963.111
$ws.Range("N67").Value = -4782.6978
# Row 76
$ws.Range("H76").Value = 3336
$ws.Range("I76").Value = 2933.3333
$ws.Range("J76").Value = 3940
$ws.Range("K76").Value = 2933.3333
$ws.Range("L76").Value = 3940
$ws.Range("M76").Value = -2618.3333
$ws.Range("N76").Value = -4570
# Row 79
$ws.Range("H79").Value = 3336
$ws.Range("I79").Value = 2933.3333
$ws.Range("J79").Value = 3940
$ws.Range("K79").Value = 2933.3333
$ws.Range("L79").Value = 3940
$ws.Range("M79").Value = -1841.3333
$ws.Range("N79").Value = -6124
# Row 138
$ws.Range("H138").Value = 3202.875
$ws.Range("I138").Value = 2072.7856
$ws.Range("J138").Value = 4785
$ws.Range("K138").Value = 6218.3568
$ws.Range("L138").Value = 14355
$ws.Range("M138").Value = -1078.3568
$ws.Range("N138").Value = -24635

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 1003795.8
$ws.Range("I63").Value = 1252245
$ws.Range("J63").Value = 9999
$ws.Range("K63").Value = 1252245
$ws.Range("L63").Value = 9999
$ws.Range("M63").Value = -1251559
$ws.Range("N63").Value = -11371
# Row 66
$ws.Range("H66").Value = 1003795.8
$ws.Range("I66").Value = 1252245
$ws.Range("J66").Value = 9999
$ws.Range("K66").Value = 6261225
$ws.Range("L66").Value = 49995
$ws.Range("M66").Value = -6257793
$ws.Range("N66").Value = -56859
# Row 122
$ws.Range("H122").Value = 2313
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2313
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6939
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -11839

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 88
$ws.Range("H88").Value = 9800
$ws.Range("J88").Value = 9800
$ws.Range("L88").Value = 9800
$ws.Range("N88").Value = -10612
# Row 91
$ws.Range("H91").Value = 9800
$ws.Range("J91").Value = 9800
$ws.Range("L91").Value = 9800
$ws.Range("N91").Value = -12608
# Row 105
$ws.Range("H105").Value = 2287.4666
$ws.Range("I105").Value = 2210.9092
$ws.Range("J105").Value = 2498
$ws.Range("K105").Value = 2210.9092
$ws.Range("L105").Value = 2498
$ws.Range("M105").Value = -463.9092000000001
$ws.Range("N105").Value = -5992

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 18817.5
$ws.Range("I62").Value = 2581
$ws.Range("J62").Value = 100000
$ws.Range("K62").Value = 2581
$ws.Range("L62").Value = 100000
$ws.Range("M62").Value = -1957
$ws.Range("N62").Value = -101248
# Row 65
$ws.Range("H65").Value = 18817.5
$ws.Range("I65").Value = 2581
$ws.Range("J65").Value = 100000
$ws.Range("K65").Value = 12905
$ws.Range("L65").Value = 500000
$ws.Range("M65").Value = -9785
$ws.Range("N65").Value = -506240
# Row 88
$ws.Range("H88").Value = 25997.777
$ws.Range("J88").Value = 27372.5
$ws.Range("L88").Value = 27372.5
$ws.Range("N88").Value = -28184.5
# Row 91
$ws.Range("H91").Value = 25997.777
$ws.Range("J91").Value = 27372.5
$ws.Range("L91").Value = 27372.5
$ws.Range("N91").Value = -30180.5
# Row 112
$ws.Range("H112").Value = 39500
$ws.Range("J112").Value = 39500
$ws.Range("L112").Value = 39500
$ws.Range("N112").Value = -42454

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 19
$ws.Range("H19").Value = 9860.5
$ws.Range("I19").Value = 720
$ws.Range("J19").Value = 19001
$ws.Range("K19").Value = 2160
$ws.Range("L19").Value = 57003
$ws.Range("M19").Value = -1986
$ws.Range("N19").Value = -57351

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
# Row 70
$ws.Range("H70").Value = 6088.054
$ws.Range("I70").Value = 6871.25
$ws.Range("J70").Value = 5491.3335
$ws.Range("K70").Value = 6871.25
$ws.Range("L70").Value = 5491.3335
$ws.Range("M70").Value = -6601.25
$ws.Range("N70").Value = -6031.3335
# Row 73
$ws.Range("H73").Value = 6088.054
$ws.Range("I73").Value = 6871.25
$ws.Range("J73").Value = 5491.3335
$ws.Range("K73").Value = 6871.25
$ws.Range("L73").Value = 5491.3335
$ws.Range("M73").Value = -5935.25
$ws.Range("N73").Value = -7363.3335
# Row 80
$ws.Range("H80").Value = 2550.75
$ws.Range("I80").Value = 2250
$ws.Range("J80").Value = 2851.5
$ws.Range("K80").Value = 2250
$ws.Range("L80").Value = 2851.5
$ws.Range("M80").Value = -1252
$ws.Range("N80").Value = -4847.5
# Row 83
$ws.Range("H83").Value = 2550.75
$ws.Range("I83").Value = 2250
$ws.Range("J83").Value = 2851.5
$ws.Range("K83").Value = 11250
$ws.Range("L83").Value = 14257.5
$ws.Range("M83").Value = -6258
$ws.Range("N83").Value = -24241.5
# Row 122
$ws.Range("H122").Value = 6787.7144
$ws.Range("I122").Value = 9238
$ws.Range("J122").Value = 4950
$ws.Range("K122").Value = 27714
$ws.Range("L122").Value = 14850
$ws.Range("M122").Value = -25264
$ws.Range("N122").Value = -19750

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4825
$ws.Range("I7").Value = 4500
$ws.Range("J7").Value = 4987.5
$ws.Range("K7").Value = 4500
$ws.Range("L7").Value = 4987.5
$ws.Range("M7").Value = -4388
$ws.Range("N7").Value = -5211.5
# Row 40
$ws.Range("H40").Value = 90913816
$ws.Range("I40").Value = 111115000
$ws.Range("J40").Value = 8490
$ws.Range("K40").Value = 111115000
$ws.Range("L40").Value = 8490
$ws.Range("M40").Value = -111114864
$ws.Range("N40").Value = -8762
# Row 87
$ws.Range("H87").Value = 18000
$ws.Range("J87").Value = 18000
$ws.Range("L87").Value = 18000
$ws.Range("N87").Value = -20246
# Row 88
$ws.Range("H88").Value = 9500
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 18000
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 18000
$ws.Range("M88").Value = -572
$ws.Range("N88").Value = -18856
# Row 90
$ws.Range("H90").Value = 18000
$ws.Range("J90").Value = 18000
$ws.Range("L90").Value = 54000
$ws.Range("N90").Value = -65232
# Row 91
$ws.Range("H91").Value = 9500
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 18000
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 18000
$ws.Range("M91").Value = 482
$ws.Range("N91").Value = -20964
# Row 126
$ws.Range("H126").Value = 4825
$ws.Range("I126").Value = 4500
$ws.Range("J126").Value = 4987.5
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 14962.5
$ws.Range("M126").Value = -11030
$ws.Range("N126").Value = -19902.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 17
$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 5000
$ws.Range("K17").Value = 5000
$ws.Range("M17").Value = -4828
